$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "sathiya"
$ws.Range("D3").Value = "sathiya@gmail.com"
$ws.Range("E3").Value = "globalTiger"

# Mobile number must stay a text value (not be coerced to a number),
# matching the original cell's string type.
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "55667788"
